$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing data rows (2-23, columns A-R / 1-18) before overwriting
# anything, so the subsequent permutation writes never read already-moved data.
$rows = @{}
for ($r = 2; $r -le 23; $r++) {
    $rowvals = @()
    for ($c = 1; $c -le 18; $c++) {
        $rowvals += ,$ws.Cells.Item($r, $c).Value2
    }
    $rows[$r] = $rowvals
}

# Target row -> source row mapping (the row reorder described by the diff).
$mapping = @{
    2  = 20
    3  = 17
    4  = 5
    5  = 8
    6  = 16
    7  = 4
    8  = 12
    9  = 13
    10 = 7
    11 = 14
    12 = 18
    13 = 15
    14 = 19
    15 = 21
    16 = 10
    17 = 6
    18 = 11
    19 = 23
    20 = 22
    21 = 3
    22 = 2
    23 = 9
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $rows[$srcRow]
    for ($c = 1; $c -le 18; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c - 1]
    }
}
